$wb = $excel.ActiveWorkbook
$zdock = $wb.Worksheets.Item("ZDock")

# Add a new worksheet named "ClusPro" right after "ZDock"
$cluspro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $zdock)
$cluspro.Name = "ClusPro"

# Copy the header row (Target, Complex 1..Complex 10) from ZDock into ClusPro
$zdock.Range("A1:K1").Copy() | Out-Null
$cluspro.Range("A1").PasteSpecial(-4163) | Out-Null

# Copy the target name column (A2:A36) from ZDock into ClusPro
$zdock.Range("A2:A36").Copy() | Out-Null
$cluspro.Range("A2").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# Update selections to match the saved state
$zdock.Range("B34").Select() | Out-Null
$cluspro.Range("B2").Select() | Out-Null

# Make sure ZDock remains the active sheet/tab
$zdock.Activate()
$zdock.Range("B34").Select() | Out-Null
